# Auto-generated edit script: refreshes market-price / profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 728.7778
$ws.Range("I28").Value = 383.75
$ws.Range("K28").Value = 383.75
$ws.Range("M28").Value = 101.25
$ws.Range("H86").Value = 6442.7
$ws.Range("I86").Value = 5158.1113
$ws.Range("J86").Value = 18004
$ws.Range("K86").Value = 5158.1113
$ws.Range("L86").Value = 18004
$ws.Range("M86").Value = -4035.1113
$ws.Range("N86").Value = -20250
$ws.Range("H89").Value = 6442.7
$ws.Range("I89").Value = 5158.1113
$ws.Range("J89").Value = 18004
$ws.Range("K89").Value = 25790.5565
$ws.Range("L89").Value = 90020
$ws.Range("M89").Value = -20174.5565
$ws.Range("N89").Value = -101252
$ws.Range("H98").Value = 1874.5834
$ws.Range("I98").Value = 1849.5
$ws.Range("K98").Value = 1849.5
$ws.Range("M98").Value = -351.5
$ws.Range("H116").Value = 5429.593
$ws.Range("I116").Value = 8213.267
$ws.Range("J116").Value = 1950
$ws.Range("K116").Value = 8213.267
$ws.Range("L116").Value = 1950
$ws.Range("M116").Value = -4771.267
$ws.Range("N116").Value = -8834
$ws.Range("H122").Value = 1874.5834
$ws.Range("I122").Value = 1849.5
$ws.Range("K122").Value = 5548.5
$ws.Range("M122").Value = -3098.5
$ws.Range("H123").Value = 54552
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 54552
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 54552
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -64352
$ws.Range("H132").Value = 2141.0454
$ws.Range("I132").Value = 2141.0454
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6423.1362
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3893.1362
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 1389.6154
$ws.Range("I137").Value = 795.8095
$ws.Range("J137").Value = 1791.871
$ws.Range("K137").Value = 2387.4285
$ws.Range("L137").Value = 5375.613
$ws.Range("M137").Value = 162.5715
$ws.Range("N137").Value = -10475.613
$ws.Range("H138").Value = 3474.047
$ws.Range("I138").Value = 1349.4375
$ws.Range("J138").Value = 4756.83
$ws.Range("K138").Value = 4048.3125
$ws.Range("L138").Value = 14270.49
$ws.Range("M138").Value = 1091.6875
$ws.Range("N138").Value = -24550.49

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7306.094
$ws.Range("I32").Value = 7769.4443
$ws.Range("J32").Value = 6324.8823
$ws.Range("K32").Value = 7769.4443
$ws.Range("L32").Value = 6324.8823
$ws.Range("M32").Value = -7482.4443
$ws.Range("N32").Value = -6898.8823
$ws.Range("H45").Value = 9721.333000000001
$ws.Range("I45").Value = 20788.4
$ws.Range("K45").Value = 20788.4
$ws.Range("M45").Value = -20411.4
$ws.Range("H61").Value = 3806.9768
$ws.Range("I61").Value = 4939.1333
$ws.Range("J61").Value = 1194.3077
$ws.Range("K61").Value = 4939.1333
$ws.Range("L61").Value = 1194.3077
$ws.Range("M61").Value = -4727.1333
$ws.Range("N61").Value = -1618.3077
$ws.Range("H109").Value = 39000
$ws.Range("J109").Value = 39000
$ws.Range("L109").Value = 39000
$ws.Range("N109").Value = -41774
$ws.Range("H122").Value = 887098.5600000001
$ws.Range("I122").Value = 988939.6
$ws.Range("J122").Value = 4476
$ws.Range("K122").Value = 2966818.8
$ws.Range("L122").Value = 13428
$ws.Range("M122").Value = -2964368.8
$ws.Range("N122").Value = -18328
$ws.Range("H132").Value = 2333.9822
$ws.Range("I132").Value = 1370.6578
$ws.Range("J132").Value = 4367.6665
$ws.Range("K132").Value = 4111.9734
$ws.Range("L132").Value = 13102.9995
$ws.Range("M132").Value = -1581.9734
$ws.Range("N132").Value = -18162.9995
$ws.Range("H136").Value = 3806.9768
$ws.Range("I136").Value = 4939.1333
$ws.Range("J136").Value = 1194.3077
$ws.Range("K136").Value = 14817.3999
$ws.Range("L136").Value = 3582.9231
$ws.Range("M136").Value = -12267.3999
$ws.Range("N136").Value = -8682.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 71430280
$ws.Range("I99").Value = 100001304
$ws.Range("J99").Value = 2725
$ws.Range("K99").Value = 100001304
$ws.Range("L99").Value = 2725
$ws.Range("M99").Value = -99999806
$ws.Range("N99").Value = -5721
$ws.Range("H107").Value = 1661.3
$ws.Range("I107").Value = 1742.8572
$ws.Range("J107").Value = 1471
$ws.Range("K107").Value = 1742.8572
$ws.Range("L107").Value = 1471
$ws.Range("M107").Value = 177.1428000000001
$ws.Range("N107").Value = -5311

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 25007202
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 1500
$ws.Range("M99").Value = -2
$ws.Range("H107").Value = 585.3823
$ws.Range("I107").Value = 493.83334
$ws.Range("J107").Value = 688.375
$ws.Range("K107").Value = 493.83334
$ws.Range("L107").Value = 688.375
$ws.Range("M107").Value = 1426.16666
$ws.Range("N107").Value = -4528.375
$ws.Range("H126").Value = 25007202
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H134").Value = 1582.375
$ws.Range("I134").Value = 1582.375
$ws.Range("K134").Value = 4747.125
$ws.Range("M134").Value = -2212.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 588758.9
$ws.Range("I8").Value = 588758.9
$ws.Range("K8").Value = 1766276.7
$ws.Range("M8").Value = -1766137.7
$ws.Range("H12").Value = 5263258.5
$ws.Range("I12").Value = 8333428
$ws.Range("J12").Value = 111.14286
$ws.Range("K12").Value = 25000284
$ws.Range("L12").Value = 333.42858
$ws.Range("M12").Value = -25000111
$ws.Range("N12").Value = -679.42858
$ws.Range("H68").Value = 2306.5745
$ws.Range("I68").Value = 2570.4443
$ws.Range("J68").Value = 1950.35
$ws.Range("K68").Value = 7711.3329
$ws.Range("L68").Value = 5851.049999999999
$ws.Range("M68").Value = -6900.3329
$ws.Range("N68").Value = -7473.049999999999
$ws.Range("H71").Value = 2306.5745
$ws.Range("I71").Value = 2570.4443
$ws.Range("J71").Value = 1950.35
$ws.Range("K71").Value = 23133.9987
$ws.Range("L71").Value = 17553.15
$ws.Range("M71").Value = -19077.9987
$ws.Range("N71").Value = -25665.15
$ws.Range("H113").Value = 1053132.1
$ws.Range("I113").Value = 1471048.1
$ws.Range("K113").Value = 4413144.300000001
$ws.Range("M113").Value = -4410974.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6089.7334
$ws.Range("I70").Value = 6090.591
$ws.Range("J70").Value = 6087.375
$ws.Range("K70").Value = 6090.591
$ws.Range("L70").Value = 6087.375
$ws.Range("M70").Value = -5820.591
$ws.Range("N70").Value = -6627.375
$ws.Range("H73").Value = 6089.7334
$ws.Range("I73").Value = 6090.591
$ws.Range("J73").Value = 6087.375
$ws.Range("K73").Value = 6090.591
$ws.Range("L73").Value = 6087.375
$ws.Range("M73").Value = -5154.591
$ws.Range("N73").Value = -7959.375
$ws.Range("H102").Value = 2585.1875
$ws.Range("I102").Value = 1112.9166
$ws.Range("J102").Value = 7002
$ws.Range("K102").Value = 1112.9166
$ws.Range("L102").Value = 7002
$ws.Range("M102").Value = 509.0834
$ws.Range("N102").Value = -10246
$ws.Range("H122").Value = 4987576
$ws.Range("I122").Value = 6482798.5
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 19448395.5
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -19445945.5
$ws.Range("N122").Value = -15400
$ws.Range("H126").Value = 6162.4346
$ws.Range("I126").Value = 6492.143
$ws.Range("J126").Value = 2700.5
$ws.Range("K126").Value = 19476.429
$ws.Range("L126").Value = 8101.5
$ws.Range("M126").Value = -17006.429
$ws.Range("N126").Value = -13041.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 338000
$ws.Range("J18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("N18").Value = -7344
$ws.Range("H20").Value = 10000000
$ws.Range("J20").Value = 10000000
$ws.Range("L20").Value = 10000000
$ws.Range("N20").Value = -10000452
$ws.Range("H122").Value = 4075383.2
$ws.Range("I122").Value = 5105527
$ws.Range("J122").Value = 1671714
$ws.Range("K122").Value = 15316581
$ws.Range("L122").Value = 5015142
$ws.Range("M122").Value = -15314131
$ws.Range("N122").Value = -5020042
$ws.Range("H128").Value = 32000
$ws.Range("J128").Value = 32000
$ws.Range("L128").Value = 32000
$ws.Range("N128").Value = -41960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2181.8948
$ws.Range("J96").Value = 2626.5
$ws.Range("L96").Value = 2626.5
$ws.Range("N96").Value = -5372.5
$ws.Range("H122").Value = 3683.2778
$ws.Range("I122").Value = 2799.9375
$ws.Range("K122").Value = 8399.8125
$ws.Range("M122").Value = -5949.8125
